# Generate Report for Handback
# This script refreshes the localization-status workbook to reflect that the
# zh-cn and de-de handback files are now in sync with en-US (i.e. a handback
# just completed successfully): the "Ready for handoff" status becomes
# "Handed back: in sync with en-US", the per-language "Latest Handback
# DateTime" timestamps move forward, and the stale "Error Detail" messages
# (about the handback file version lagging behind the latest source) are
# cleared out now that everything is current. Some report columns are also
# resized to better fit the (now shorter/longer) content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for both rows ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn detail sheet ---
# Status column (C) reflects the same refreshed status text.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
# Latest Handback DateTime (K) advances to the new handback timestamp.
$wsZhCn.Range("K2").Value = "2016-08-05 02:33:42"
$wsZhCn.Range("K3").Value = "2016-08-05 02:33:42"
# Error Detail (P) is cleared now that the handback file is up to date.
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

# --- de-de detail sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-05 02:34:00"
$wsDeDe.Range("K3").Value = "2016-08-05 02:34:00"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

# --- Column width adjustments (report columns resized to fit content) ---
# Overview: zh-cn / de-de status columns widened.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14437166849777
$wsOverview.Columns.Item(6).ColumnWidth = 29.14437166849777

# zh-cn / de-de: Status column (C) widened, Error Detail column (P) narrowed
# now that it no longer holds long error text.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719813028965

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719813028965
